# Updated cryptos list — refresh Price (col D) and Volume(1h) (col E) figures,
# and swap the TrustWalletToken / InternetComputer(DFINITY) rows (37/38).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '24.681.54'; ForceText = $false },
    @{ Cell = 'E2'; Value = '  -1.05%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '1.658.74'; ForceText = $false },
    @{ Cell = 'E3'; Value = '  -2.90%  '; ForceText = $false },
    @{ Cell = 'D4'; Value = '1.001'; ForceText = $true },
    @{ Cell = 'E4'; Value = '  -0.28%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '319.98'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +2.19%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '0.9985'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  +0.06%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '0.3636'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -2.92%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '47.14'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -4.73%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '0.3270'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  -5.00%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '1.133'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  -7.59%  '; ForceText = $false },
    @{ Cell = 'D11'; Value = '0.07055'; ForceText = $true },
    @{ Cell = 'E11'; Value = '  -6.42%  '; ForceText = $false },
    @{ Cell = 'D12'; Value = '0.9978'; ForceText = $true },
    @{ Cell = 'D13'; Value = '5.983'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  -5.40%  '; ForceText = $false },
    @{ Cell = 'D14'; Value = '19.55'; ForceText = $true },
    @{ Cell = 'E14'; Value = '  -8.01%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '1.660.87'; ForceText = $false },
    @{ Cell = 'D16'; Value = '6.627'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  -6.37%  '; ForceText = $false },
    @{ Cell = 'E17'; Value = '  -7.48%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '0.06589'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  -2.10%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '0.9975'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  +0.04%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '78.88'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  -6.26%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '5.947'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  -6.99%  '; ForceText = $false },
    @{ Cell = 'D22'; Value = '15.78'; ForceText = $true },
    @{ Cell = 'E22'; Value = '  -8.91%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '12.62'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  -3.49%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '24.646.58'; ForceText = $false },
    @{ Cell = 'E24'; Value = '  -1.15%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '2.465'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +0.98%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '2.413'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  -13.79%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '148.14'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -1.17%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '18.61'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  -8.87%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '1.845.02'; ForceText = $false },
    @{ Cell = 'E29'; Value = '  -2.61%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '1.218'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  -2.80%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '125.15'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  -5.86%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '4.072'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  -3.59%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '5.840'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  -14.53%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '0.08472'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  -3.84%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '1.681'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  -5.15%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '12.38'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  -10.95%  '; ForceText = $false },
    @{ Cell = 'B37'; Value = 'TrustWalletToken'; ForceText = $false },
    @{ Cell = 'C37'; Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'; ForceText = $false },
    @{ Cell = 'D37'; Value = '1.275'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  +0.29%  '; ForceText = $false },
    @{ Cell = 'B38'; Value = 'InternetComputer(DFINITY)'; ForceText = $false },
    @{ Cell = 'C38'; Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'; ForceText = $false },
    @{ Cell = 'D38'; Value = '5.219'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  -7.42%  '; ForceText = $false },
    @{ Cell = 'D39'; Value = '0.06052'; ForceText = $true },
    @{ Cell = 'E39'; Value = '  -9.27%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '0.02240'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  -7.29%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '0.2076'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  -7.25%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '8.203'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  -10.70%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '0.9975'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +0.00%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '0.5930'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  -8.50%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '3.851'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  +0.19%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '12.74'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  -8.12%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '0.5619'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  -8.92%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '124.62'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  -3.72%  '; ForceText = $false },
    @{ Cell = 'D49'; Value = '1.957'; ForceText = $true },
    @{ Cell = 'E49'; Value = '  -7.98%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '0.06979'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  -4.83%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '1.194'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  -3.90%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.ForceText) {
        # Several "Price" entries look numeric (e.g. "1.001", "0.9985") but are
        # stored as literal text in the source sheet -- force text formatting
        # before writing so Excel doesn't silently coerce them to numbers.
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
